$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3463.4849
$ws.Range("I137").Value = 2077.4546
$ws.Range("K137").Value = 6232.3638
$ws.Range("M137").Value = -3682.3638
$ws.Range("H138").Value = 3529.53
$ws.Range("I138").Value = 676
$ws.Range("J138").Value = 5431.8833
$ws.Range("K138").Value = 2028
$ws.Range("L138").Value = 16295.6499
$ws.Range("M138").Value = 3112
$ws.Range("N138").Value = -26575.6499

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3515.2703
$ws.Range("I32").Value = 3159.682
$ws.Range("J32").Value = 6448.875
$ws.Range("K32").Value = 3159.682
$ws.Range("L32").Value = 6448.875
$ws.Range("M32").Value = -2872.682
$ws.Range("N32").Value = -7022.875
$ws.Range("H61").Value = 3050
$ws.Range("I61").Value = 1800
$ws.Range("J61").Value = 3466.6667
$ws.Range("K61").Value = 1800
$ws.Range("L61").Value = 3466.6667
$ws.Range("M61").Value = -1588
$ws.Range("N61").Value = -3890.6667
$ws.Range("H63").Value = 7293109
$ws.Range("I63").Value = 10656713
$ws.Range("K63").Value = 10656713
$ws.Range("M63").Value = -10656027
$ws.Range("H66").Value = 7293109
$ws.Range("I66").Value = 10656713
$ws.Range("K66").Value = 53283565
$ws.Range("M66").Value = -53280133
$ws.Range("H74").Value = 4210.7666
$ws.Range("I74").Value = 4396.304
$ws.Range("J74").Value = 3601.1428
$ws.Range("K74").Value = 4396.304
$ws.Range("L74").Value = 3601.1428
$ws.Range("M74").Value = -3522.304
$ws.Range("N74").Value = -5349.1428
$ws.Range("H77").Value = 4210.7666
$ws.Range("I77").Value = 4396.304
$ws.Range("J77").Value = 3601.1428
$ws.Range("K77").Value = 21981.52
$ws.Range("L77").Value = 18005.714
$ws.Range("M77").Value = -17613.52
$ws.Range("N77").Value = -26741.714
$ws.Range("H88").Value = 11113390
$ws.Range("I88").Value = 33334232
$ws.Range("J88").Value = 2969.75
$ws.Range("K88").Value = 33334232
$ws.Range("L88").Value = 2969.75
$ws.Range("M88").Value = -33333826
$ws.Range("N88").Value = -3781.75
$ws.Range("H91").Value = 11113390
$ws.Range("I91").Value = 33334232
$ws.Range("J91").Value = 2969.75
$ws.Range("K91").Value = 33334232
$ws.Range("L91").Value = 2969.75
$ws.Range("M91").Value = -33332828
$ws.Range("N91").Value = -5777.75
$ws.Range("H132").Value = 2263.5098
$ws.Range("I132").Value = 1374.4324
$ws.Range("J132").Value = 4613.2144
$ws.Range("K132").Value = 4123.2972
$ws.Range("L132").Value = 13839.6432
$ws.Range("M132").Value = -1593.2972
$ws.Range("N132").Value = -18899.6432
$ws.Range("H136").Value = 3050
$ws.Range("I136").Value = 1800
$ws.Range("J136").Value = 3466.6667
$ws.Range("K136").Value = 5400
$ws.Range("L136").Value = 10400.0001
$ws.Range("M136").Value = -2850
$ws.Range("N136").Value = -15500.0001
$ws.Range("H138").Value = 78739
$ws.Range("J138").Value = 78739
$ws.Range("L138").Value = 78739
$ws.Range("N138").Value = -89019
$ws.Range("H140").Value = 115000
$ws.Range("J140").Value = 115000
$ws.Range("L140").Value = 115000
$ws.Range("N140").Value = -125360

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1350.2
$ws.Range("I86").Value = 1159
$ws.Range("J86").Value = 1876
$ws.Range("K86").Value = 1159
$ws.Range("L86").Value = 1876
$ws.Range("M86").Value = -36
$ws.Range("N86").Value = -4122
$ws.Range("H89").Value = 1350.2
$ws.Range("I89").Value = 1159
$ws.Range("J89").Value = 1876
$ws.Range("K89").Value = 5795
$ws.Range("L89").Value = 9380
$ws.Range("M89").Value = -179
$ws.Range("N89").Value = -20612
$ws.Range("H134").Value = 1902.4667
$ws.Range("I134").Value = 1215.4
$ws.Range("J134").Value = 3276.6
$ws.Range("K134").Value = 3646.2
$ws.Range("L134").Value = 9829.799999999999
$ws.Range("M134").Value = -1111.2
$ws.Range("N134").Value = -14899.8

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 496.14706
$ws.Range("I22").Value = 308.8
$ws.Range("J22").Value = 1016.55554
$ws.Range("K22").Value = 308.8
$ws.Range("L22").Value = 1016.55554
$ws.Range("M22").Value = 41.19999999999999
$ws.Range("N22").Value = -1716.55554
$ws.Range("H31").Value = 7145276.5
$ws.Range("I31").Value = 1336.2273
$ws.Range("J31").Value = 19235022
$ws.Range("K31").Value = 1336.2273
$ws.Range("L31").Value = 19235022
$ws.Range("M31").Value = -1041.2273
$ws.Range("N31").Value = -19235612
$ws.Range("H34").Value = 7145276.5
$ws.Range("I34").Value = 1336.2273
$ws.Range("J34").Value = 19235022
$ws.Range("K34").Value = 1336.2273
$ws.Range("L34").Value = 19235022
$ws.Range("M34").Value = -1134.2273
$ws.Range("N34").Value = -19235426
$ws.Range("H58").Value = 2022.6528
$ws.Range("I58").Value = 1819.7455
$ws.Range("J58").Value = 2679.1177
$ws.Range("K58").Value = 1819.7455
$ws.Range("L58").Value = 2679.1177
$ws.Range("M58").Value = -1616.7455
$ws.Range("N58").Value = -3085.1177
$ws.Range("H105").Value = 2132
$ws.Range("I105").Value = 1957.1428
$ws.Range("K105").Value = 1957.1428
$ws.Range("M105").Value = -210.1428000000001
$ws.Range("H134").Value = 4099.023
$ws.Range("I134").Value = 5030.28
$ws.Range("J134").Value = 2873.6843
$ws.Range("K134").Value = 15090.84
$ws.Range("L134").Value = 8621.052899999999
$ws.Range("M134").Value = -12555.84
$ws.Range("N134").Value = -13691.0529
$ws.Range("H136").Value = 2022.6528
$ws.Range("I136").Value = 1819.7455
$ws.Range("J136").Value = 2679.1177
$ws.Range("K136").Value = 5459.2365
$ws.Range("L136").Value = 8037.353099999999
$ws.Range("M136").Value = -2909.2365
$ws.Range("N136").Value = -13137.3531

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1429.8966
$ws.Range("J5").Value = 3467.3
$ws.Range("L5").Value = 10401.9
$ws.Range("N5").Value = -10625.9
$ws.Range("H22").Value = 1542.8572
$ws.Range("J22").Value = 1466.6666
$ws.Range("L22").Value = 4399.9998
$ws.Range("N22").Value = -4737.9998
$ws.Range("H27").Value = 1542.8572
$ws.Range("J27").Value = 1466.6666
$ws.Range("L27").Value = 4399.9998
$ws.Range("N27").Value = -4603.9998
$ws.Range("H49").Value = 2772.3635
$ws.Range("I49").Value = 2501.5
$ws.Range("J49").Value = 2927.1428
$ws.Range("K49").Value = 7504.5
$ws.Range("L49").Value = 8781.428400000001
$ws.Range("M49").Value = -7348.5
$ws.Range("N49").Value = -9093.428400000001
$ws.Range("H107").Value = 24843.049
$ws.Range("J107").Value = 38921
$ws.Range("L107").Value = 116763
$ws.Range("N107").Value = -120603
$ws.Range("H113").Value = 551.6719000000001
$ws.Range("I113").Value = 500.48486
$ws.Range("J113").Value = 606.1613
$ws.Range("K113").Value = 1501.45458
$ws.Range("L113").Value = 1818.4839
$ws.Range("M113").Value = 668.5454199999999
$ws.Range("N113").Value = -6158.4839
$ws.Range("H131").Value = 883.7273
$ws.Range("I131").Value = 427.22223
$ws.Range("J131").Value = 955.807
$ws.Range("K131").Value = 1281.66669
$ws.Range("L131").Value = 2867.421
$ws.Range("M131").Value = 3758.33331
$ws.Range("N131").Value = -12947.421
$ws.Range("H135").Value = 1429.8966
$ws.Range("J135").Value = 3467.3
$ws.Range("L135").Value = 31205.7
$ws.Range("N135").Value = -36275.7

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6303.905
$ws.Range("I70").Value = 5743.6333
$ws.Range("K70").Value = 5743.6333
$ws.Range("M70").Value = -5473.6333
$ws.Range("H73").Value = 6303.905
$ws.Range("I73").Value = 5743.6333
$ws.Range("K73").Value = 5743.6333
$ws.Range("M73").Value = -4807.6333
$ws.Range("H122").Value = 5011.1816
$ws.Range("J122").Value = 6027
$ws.Range("L122").Value = 18081
$ws.Range("N122").Value = -22981
$ws.Range("H132").Value = 3155.1562
$ws.Range("I132").Value = 2074.1177
$ws.Range("J132").Value = 4380.3335
$ws.Range("K132").Value = 6222.353099999999
$ws.Range("L132").Value = 13141.0005
$ws.Range("M132").Value = -3692.353099999999
$ws.Range("N132").Value = -18201.0005

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 946.6875
$ws.Range("I16").Value = 795.9259
$ws.Range("J16").Value = 1760.8
$ws.Range("K16").Value = 795.9259
$ws.Range("L16").Value = 1760.8
$ws.Range("M16").Value = -625.9259
$ws.Range("N16").Value = -2100.8
$ws.Range("H61").Value = 1122.7693
$ws.Range("I61").Value = 1031.7
$ws.Range("J61").Value = 1426.3334
$ws.Range("K61").Value = 1031.7
$ws.Range("L61").Value = 1426.3334
$ws.Range("M61").Value = -829.7
$ws.Range("N61").Value = -1830.3334
$ws.Range("H113").Value = 1122.7693
$ws.Range("I113").Value = 1031.7
$ws.Range("J113").Value = 1426.3334
$ws.Range("K113").Value = 1031.7
$ws.Range("L113").Value = 1426.3334
$ws.Range("M113").Value = 1138.3
$ws.Range("N113").Value = -5766.3334
$ws.Range("H136").Value = 2475.5557
$ws.Range("I136").Value = 1422.8572
$ws.Range("K136").Value = 4268.571599999999
$ws.Range("M136").Value = -1718.571599999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 395.0625
$ws.Range("I113").Value = 276.5
$ws.Range("J113").Value = 513.625
$ws.Range("K113").Value = 829.5
$ws.Range("L113").Value = 1540.875
$ws.Range("M113").Value = 1340.5
$ws.Range("N113").Value = -5880.875
$ws.Range("H132").Value = 8131989
$ws.Range("I132").Value = 885.6667
$ws.Range("J132").Value = 19611194
$ws.Range("K132").Value = 2657.0001
$ws.Range("L132").Value = 58833582
$ws.Range("M132").Value = -127.0001000000002
$ws.Range("N132").Value = -58838642
$ws.Range("H136").Value = 1578.9125
$ws.Range("I136").Value = 562.05554
$ws.Range("K136").Value = 1686.16662
$ws.Range("M136").Value = 863.83338
